{"js": "// Replace the division-problem answers in the table with the new set of\n// problems/answers. Each \"old\" string is unique in the document, so a\n// simple exact-text search + Replace is sufficient and keeps the original\n// run formatting (font/size) untouched.\nconst replacements = [\n  [\"87\u00f77=12, 3\", \"19\u00f75=3, 4\"],\n  [\"30\u00f75=6, 0\", \"19\u00f72=9, 1\"],\n  [\"87\u00f79=9, 6\", \"95\u00f75=19, 0\"],\n  [\"18\u00f72=9, 0\", \"96\u00f72=48, 0\"],\n  [\"17\u00f74=4, 1\", \"68\u00f77=9, 5\"],\n  [\"51\u00f79=5, 6\", \"76\u00f72=38, 0\"],\n  [\"76\u00f75=15, 1\", \"23\u00f73=7, 2\"],\n  [\"52\u00f77=7, 3\", \"23\u00f72=11, 1\"],\n  [\"25\u00f75=5, 0\", \"76\u00f73=25, 1\"],\n  [\"21\u00f78=2, 5\", \"83\u00f74=20, 3\"],\n  [\"68\u00f73=22, 2\", \"48\u00f73=16, 0\"],\n  [\"86\u00f74=21, 2\", \"20\u00f79=2, 2\"],\n  [\"60\u00f79=6, 6\", \"99\u00f73=33, 0\"],\n  [\"72\u00f76=12, 0\", \"35\u00f77=5, 0\"],\n  [\"33\u00f72=16, 1\", \"49\u00f77=7, 0\"],\n  [\"94\u00f77=13, 3\", \"13\u00f76=2, 1\"],\n  [\"32\u00f73=10, 2\", \"19\u00f79=2, 1\"],\n  [\"31\u00f76=5, 1\", \"84\u00f73=28, 0\"],\n  [\"74\u00f78=9, 2\", \"91\u00f75=18, 1\"],\n  [\"27\u00f78=3, 3\", \"66\u00f77=9, 3\"],\n  [\"19\u00f77=2, 5\", \"51\u00f75=10, 1\"],\n  [\"60\u00f72=30, 0\", \"23\u00f76=3, 5\"],\n  [\"25\u00f72=12, 1\", \"47\u00f73=15, 2\"],\n  [\"14\u00f79=1, 5\", \"64\u00f78=8, 0\"],\n  [\"13\u00f72=6, 1\", \"80\u00f78=10, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem answers in the table with the new set of\n# problems/answers. Each \"old\" string is unique in the document, so a\n# simple Find/Replace (wdReplaceAll) per pair is sufficient and preserves\n# the original run formatting (font/size) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"87\u00f77=12, 3\", \"19\u00f75=3, 4\"),\n    @(\"30\u00f75=6, 0\", \"19\u00f72=9, 1\"),\n    @(\"87\u00f79=9, 6\", \"95\u00f75=19, 0\"),\n    @(\"18\u00f72=9, 0\", \"96\u00f72=48, 0\"),\n    @(\"17\u00f74=4, 1\", \"68\u00f77=9, 5\"),\n    @(\"51\u00f79=5, 6\", \"76\u00f72=38, 0\"),\n    @(\"76\u00f75=15, 1\", \"23\u00f73=7, 2\"),\n    @(\"52\u00f77=7, 3\", \"23\u00f72=11, 1\"),\n    @(\"25\u00f75=5, 0\", \"76\u00f73=25, 1\"),\n    @(\"21\u00f78=2, 5\", \"83\u00f74=20, 3\"),\n    @(\"68\u00f73=22, 2\", \"48\u00f73=16, 0\"),\n    @(\"86\u00f74=21, 2\", \"20\u00f79=2, 2\"),\n    @(\"60\u00f79=6, 6\", \"99\u00f73=33, 0\"),\n    @(\"72\u00f76=12, 0\", \"35\u00f77=5, 0\"),\n    @(\"33\u00f72=16, 1\", \"49\u00f77=7, 0\"),\n    @(\"94\u00f77=13, 3\", \"13\u00f76=2, 1\"),\n    @(\"32\u00f73=10, 2\", \"19\u00f79=2, 1\"),\n    @(\"31\u00f76=5, 1\", \"84\u00f73=28, 0\"),\n    @(\"74\u00f78=9, 2\", \"91\u00f75=18, 1\"),\n    @(\"27\u00f78=3, 3\", \"66\u00f77=9, 3\"),\n    @(\"19\u00f77=2, 5\", \"51\u00f75=10, 1\"),\n    @(\"60\u00f72=30, 0\", \"23\u00f76=3, 5\"),\n    @(\"25\u00f72=12, 1\", \"47\u00f73=15, 2\"),\n    @(\"14\u00f79=1, 5\", \"64\u00f78=8, 0\"),\n    @(\"13\u00f72=6, 1\", \"80\u00f78=10, 0\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute(\n        $oldText,   # FindText\n        $true,      # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
